$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (D) and 1h volume change (E) figures for rows 2-50.
# A handful of Price cells (D13, D14, D18, D21, D45, D46) contain values that
# Excel's automatic type detection would otherwise mangle (dropped trailing
# zeros, or collapsed into scientific notation for very small magnitudes), so
# those specific cells are explicitly forced to Text format before the value
# is written, preserving the exact displayed string.
$textFormatRows = @(13, 14, 18, 21, 45, 46)
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "25.809.02"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.740.08"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "228.75"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.5173"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("D8").Value = "0.2744"
$ws.Range("E8").Value = "  +4.46%  "
$ws.Range("D9").Value = "38.42"
$ws.Range("E9").Value = "  -6.43%  "
$ws.Range("D10").Value = "0.06102"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "1.740.99"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "0.07001"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "15.30"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "0.6330"
$ws.Range("E14").Value = "  +6.04%  "
$ws.Range("D15").Value = "4.497"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "76.42"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "25.829.25"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "0.000006594"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").Value = "1.960.95"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "4.045"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "8.454"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").Value = "5.112"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "136.59"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("D28").Value = "1.817"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "15.02"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "102.83"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "0.08326"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "3.618"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "3.377"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "0.04418"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "2.607"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "0.9735"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "0.5969"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "0.01563"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "1.942"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "0.9987"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "101.93"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "0.3808"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "0.7244"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "4.870"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "0.05480"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "6.246"
$ws.Range("E47").Value = "  +5.83%  "
$ws.Range("D48").Value = "0.1099"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").Value = "29.82"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "51.86"
$ws.Range("E50").Value = "  -0.35%  "
